$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D-column price cells keep their original text storage (not
# auto-converted to numbers) by forcing Text number format before writing.
$dCells = @("D2","D3","D5","D6","D7","D8","D9","D10","D11","D13","D14","D15","D17","D18","D19","D20","D21","D23","D28","D31","D32","D37","D39","D40","D41","D42","D43","D44","D46","D47","D48","D50","D51")
foreach ($ref in $dCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = '64.381.41'
$ws.Range("E2").Value = '  +0.81%  '

$ws.Range("D3").Value = '2.762.58'
$ws.Range("E3").Value = '  +0.49%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").Value = '577.87'
$ws.Range("E5").Value = '  +0.11%  '

$ws.Range("D6").Value = '160.56'
$ws.Range("E6").Value = '  +1.01%  '

$ws.Range("D7").Value = '0.998'
$ws.Range("E7").Value = '  -0.04%  '

$ws.Range("D8").Value = '0.603'
$ws.Range("E8").Value = '  -1.23%  '

$ws.Range("D9").Value = '0.111'
$ws.Range("E9").Value = '  -1.22%  '

$ws.Range("D10").Value = '0.167'
$ws.Range("E10").Value = '  +5.06%  '

$ws.Range("D11").Value = '5.85'
$ws.Range("E11").Value = '  +3.96%  '

$ws.Range("E12").Value = '  -1.05%  '

$ws.Range("D13").Value = '3.249.21'
$ws.Range("E13").Value = '  +0.44%  '

$ws.Range("D14").Value = '27.41'
$ws.Range("E14").Value = '  +1.75%  '

$ws.Range("D15").Value = '64.017.18'
$ws.Range("E15").Value = '  +0.39%  '

$ws.Range("E16").Value = '  -1.81%  '

$ws.Range("D17").Value = '2.767.34'
$ws.Range("E17").Value = '  +0.61%  '

$ws.Range("D18").Value = '12.18'
$ws.Range("E18").Value = '  -0.37%  '

$ws.Range("D19").Value = '4.86'
$ws.Range("E19").Value = '  -1.78%  '

$ws.Range("D20").Value = '358.81'
$ws.Range("E20").Value = '  -0.64%  '

$ws.Range("D21").Value = '6.67'
$ws.Range("E21").Value = '  -3.11%  '

$ws.Range("D23").Value = '0.530'
$ws.Range("E23").Value = '  -6.14%  '

$ws.Range("E24").Value = '  -1.88%  '

$ws.Range("E25").Value = '  -0.92%  '

$ws.Range("E26").Value = '  -0.39%  '

$ws.Range("E27").Value = '  -0.13%  '

$ws.Range("D28").Value = '0.0₃0925'
$ws.Range("E28").Value = '  -1.03%  '

$ws.Range("E29").Value = '  +3.62%  '

$ws.Range("E30").Value = '  -0.93%  '

$ws.Range("D31").Value = '1.37'
$ws.Range("E31").Value = '  +9.39%  '

$ws.Range("D32").Value = '168.39'
$ws.Range("E32").Value = '  +0.18%  '

$ws.Range("E33").Value = '  +3.35%  '

$ws.Range("E34").Value = '  -0.30%  '

$ws.Range("E35").Value = '  -1.62%  '

$ws.Range("E36").Value = '  -0.02%  '

$ws.Range("D37").Value = '1.84'
$ws.Range("E37").Value = '  +1.14%  '

$ws.Range("E38").Value = '  -0.88%  '

$ws.Range("D39").Value = '352.01'
$ws.Range("E39").Value = '  +6.08%  '

$ws.Range("D40").Value = '6.43'
$ws.Range("E40").Value = '  +4.39%  '

$ws.Range("D41").Value = '4.20'
$ws.Range("E41").Value = '  -0.26%  '

$ws.Range("D42").Value = '38.99'
$ws.Range("E42").Value = '  -1.41%  '

$ws.Range("D43").Value = '22.54'
$ws.Range("E43").Value = '  +2.77%  '

$ws.Range("D44").Value = '21.57'
$ws.Range("E44").Value = '  -2.11%  '

$ws.Range("E45").Value = '  -0.45%  '

$ws.Range("D46").Value = '136.86'
$ws.Range("E46").Value = '  +0.01%  '

$ws.Range("D47").Value = '0.631'
$ws.Range("E47").Value = '  -1.49%  '

$ws.Range("D48").Value = '0.0254'
$ws.Range("E48").Value = '  -1.70%  '

$ws.Range("E49").Value = '  -0.88%  '

$ws.Range("B50").Value = 'Maker'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D50").Value = '2.151.30'
$ws.Range("E50").Value = '  +1.40%  '

$ws.Range("B51").Value = 'FirstDigitalUSD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D51").Value = '0.998'
$ws.Range("E51").Value = '  -0.24%  '
